$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 4-10 to match repulled data
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 1
